# Auto-generated edit script applying numeric updates from the commit diff
# to Sheets/Seraph_Profits.xlsx (workbook tabs: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2706.373
$ws.Range("J17").Value = 2706.373
$ws.Range("L17").Value = 8119.119000000001
$ws.Range("N17").Value = -8455.119000000001

$ws.Range("H80").Value = 224.27272
$ws.Range("I80").Value = 269.3
$ws.Range("J80").Value = 186.75
$ws.Range("K80").Value = 807.9000000000001
$ws.Range("L80").Value = 560.25
$ws.Range("M80").Value = 190.0999999999999
$ws.Range("N80").Value = -2556.25

$ws.Range("H83").Value = 224.27272
$ws.Range("I83").Value = 269.3
$ws.Range("J83").Value = 186.75
$ws.Range("K83").Value = 2423.7
$ws.Range("L83").Value = 1680.75
$ws.Range("M83").Value = 2568.3
$ws.Range("N83").Value = -11664.75

$ws.Range("H112").Value = 2155.2222
$ws.Range("J112").Value = 2155.2222
$ws.Range("L112").Value = 6465.6666
$ws.Range("N112").Value = -8681.6666

$ws.Range("H116").Value = 8750
$ws.Range("I116").Value = 7500
$ws.Range("K116").Value = 7500
$ws.Range("M116").Value = -4058

$ws.Range("H132").Value = 1550.1708
$ws.Range("I132").Value = 1584.8889
$ws.Range("K132").Value = 4754.6667
$ws.Range("M132").Value = -2224.6667

$ws.Range("H137").Value = 2108.3704
$ws.Range("I137").Value = 1831.8572
$ws.Range("J137").Value = 2406.1538
$ws.Range("K137").Value = 5495.571599999999
$ws.Range("L137").Value = 7218.4614
$ws.Range("M137").Value = -2945.571599999999
$ws.Range("N137").Value = -12318.4614

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1912.1364
$ws.Range("I61").Value = 1742.8889
$ws.Range("K61").Value = 1742.8889
$ws.Range("M61").Value = -1530.8889

$ws.Range("H74").Value = 1472.4
$ws.Range("I74").Value = 956.8148
$ws.Range("J74").Value = 6112.6665
$ws.Range("K74").Value = 956.8148
$ws.Range("L74").Value = 6112.6665
$ws.Range("M74").Value = -82.81479999999999
$ws.Range("N74").Value = -7860.6665

$ws.Range("H77").Value = 1472.4
$ws.Range("I77").Value = 956.8148
$ws.Range("J77").Value = 6112.6665
$ws.Range("K77").Value = 4784.074
$ws.Range("L77").Value = 30563.3325
$ws.Range("M77").Value = -416.0739999999996
$ws.Range("N77").Value = -39299.3325

$ws.Range("H88").Value = 1338.3636
$ws.Range("I88").Value = 1200
$ws.Range("J88").Value = 1417.4286
$ws.Range("K88").Value = 1200
$ws.Range("L88").Value = 1417.4286
$ws.Range("M88").Value = -794
$ws.Range("N88").Value = -2229.4286

$ws.Range("H91").Value = 1338.3636
$ws.Range("I91").Value = 1200
$ws.Range("J91").Value = 1417.4286
$ws.Range("K91").Value = 1200
$ws.Range("L91").Value = 1417.4286
$ws.Range("M91").Value = 204
$ws.Range("N91").Value = -4225.4286

$ws.Range("H102").Value = 221.5
$ws.Range("I102").Value = 221.5
$ws.Range("K102").Value = 221.5
$ws.Range("M102").Value = 1400.5

$ws.Range("H122").Value = 458019.2
$ws.Range("I122").Value = 772378.6
$ws.Range("K122").Value = 2317135.8
$ws.Range("M122").Value = -2314685.8

$ws.Range("H132").Value = 2090.5
$ws.Range("I132").Value = 1943.6154
$ws.Range("K132").Value = 5830.8462
$ws.Range("M132").Value = -3300.8462

$ws.Range("H136").Value = 1912.1364
$ws.Range("I136").Value = 1742.8889
$ws.Range("K136").Value = 5228.6667
$ws.Range("M136").Value = -2678.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1215.2354
$ws.Range("I86").Value = 1215.2354
$ws.Range("K86").Value = 1215.2354
$ws.Range("M86").Value = -92.23540000000003

$ws.Range("H89").Value = 1215.2354
$ws.Range("I89").Value = 1215.2354
$ws.Range("K89").Value = 6076.177
$ws.Range("M89").Value = -460.1769999999997

$ws.Range("H128").Value = 3999
$ws.Range("I128").Value = 3999
$ws.Range("K128").Value = 11997
$ws.Range("M128").Value = -9507

$ws.Range("H134").Value = 3390.5386
$ws.Range("I134").Value = 3200.125
$ws.Range("K134").Value = 9600.375
$ws.Range("M134").Value = -7065.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5726.892
$ws.Range("I31").Value = 4275.3335
$ws.Range("J31").Value = 6423.64
$ws.Range("K31").Value = 4275.3335
$ws.Range("L31").Value = 6423.64
$ws.Range("M31").Value = -3980.3335
$ws.Range("N31").Value = -7013.64

$ws.Range("H34").Value = 5726.892
$ws.Range("I34").Value = 4275.3335
$ws.Range("J34").Value = 6423.64
$ws.Range("K34").Value = 4275.3335
$ws.Range("L34").Value = 6423.64
$ws.Range("M34").Value = -4073.3335
$ws.Range("N34").Value = -6827.64

$ws.Range("H41").Value = 792
$ws.Range("I41").Value = 792
$ws.Range("K41").Value = 792
$ws.Range("M41").Value = -364

$ws.Range("H47").Value = 64
$ws.Range("I47").Value = 64
$ws.Range("K47").Value = 64
$ws.Range("M47").Value = 502

$ws.Range("H58").Value = 3194.8
$ws.Range("I58").Value = 1125
$ws.Range("K58").Value = 1125
$ws.Range("M58").Value = -922

$ws.Range("H99").Value = 11973.087
$ws.Range("I99").Value = 7114
$ws.Range("K99").Value = 7114
$ws.Range("M99").Value = -5616

$ws.Range("H122").Value = 2861.4375
$ws.Range("I122").Value = 2991.1538
$ws.Range("J122").Value = 2299.3333
$ws.Range("K122").Value = 8973.4614
$ws.Range("L122").Value = 6897.999899999999
$ws.Range("M122").Value = -6523.4614
$ws.Range("N122").Value = -11797.9999

$ws.Range("H126").Value = 11973.087
$ws.Range("I126").Value = 7114
$ws.Range("K126").Value = 21342
$ws.Range("M126").Value = -18872

$ws.Range("H132").Value = 2591.4849
$ws.Range("I132").Value = 2113.6453
$ws.Range("K132").Value = 6340.9359
$ws.Range("M132").Value = -3810.9359

$ws.Range("H134").Value = 3438.4666
$ws.Range("I134").Value = 2714.4285
$ws.Range("K134").Value = 8143.2855
$ws.Range("M134").Value = -5608.2855

$ws.Range("H136").Value = 3194.8
$ws.Range("I136").Value = 1125
$ws.Range("K136").Value = 3375
$ws.Range("M136").Value = -825

$ws.Range("H141").Value = 134468.5
$ws.Range("J141").Value = 149692
$ws.Range("L141").Value = 149692
$ws.Range("N141").Value = -160052

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 24535224
$ws.Range("I4").Value = 30663564
$ws.Range("K4").Value = 91990692
$ws.Range("M4").Value = -91990580

$ws.Range("H7").Value = 12500109
$ws.Range("J7").Value = 125
$ws.Range("L7").Value = 375
$ws.Range("N7").Value = -599

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 429.29413
$ws.Range("J2").Value = 633.8182
$ws.Range("L2").Value = 633.8182
$ws.Range("N2").Value = -859.8182

$ws.Range("H11").Value = 1284751
$ws.Range("I11").Value = 1693000
$ws.Range("J11").Value = 60004
$ws.Range("K11").Value = 1693000
$ws.Range("L11").Value = 60004
$ws.Range("M11").Value = -1692861
$ws.Range("N11").Value = -60282

$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()

$ws.Range("H132").Value = 3127
$ws.Range("I132").Value = 2285
$ws.Range("K132").Value = 6855
$ws.Range("M132").Value = -4325

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6000
$ws.Range("I7").Value = 5000
$ws.Range("J7").Value = 7000
$ws.Range("K7").Value = 5000
$ws.Range("L7").Value = 7000
$ws.Range("M7").Value = -4888
$ws.Range("N7").Value = -7224

$ws.Range("H34").Value = 4933.3335
$ws.Range("I34").Value = 4933.3335
$ws.Range("K34").Value = 4933.3335
$ws.Range("M34").Value = -4761.3335

$ws.Range("H46").Value = 3837.25
$ws.Range("I46").Value = 3000
$ws.Range("K46").Value = 3000
$ws.Range("M46").Value = -2812

$ws.Range("H100").Value = 1564.1428
$ws.Range("I100").Value = 999.6667
$ws.Range("K100").Value = 999.6667
$ws.Range("M100").Value = -458.6667

$ws.Range("H126").Value = 6000
$ws.Range("I126").Value = 5000
$ws.Range("J126").Value = 7000
$ws.Range("K126").Value = 15000
$ws.Range("L126").Value = 21000
$ws.Range("M126").Value = -12530
$ws.Range("N126").Value = -25940

$ws.Range("H132").Value = 4201.357
$ws.Range("I132").Value = 3693.647
$ws.Range("J132").Value = 4986
$ws.Range("K132").Value = 11080.941
$ws.Range("L132").Value = 14958
$ws.Range("M132").Value = -8550.940999999999
$ws.Range("N132").Value = -20018

$ws.Range("H136").Value = 1002
$ws.Range("I136").Value = 1002
$ws.Range("K136").Value = 3006
$ws.Range("M136").Value = -456

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1005.36365
$ws.Range("I113").Value = 955.9
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 2867.7
$ws.Range("L113").Value = 4500
$ws.Range("M113").Value = -697.6999999999998
$ws.Range("N113").Value = -8840

$ws.Range("H132").Value = 1125.5625
$ws.Range("I132").Value = 1125.5625
$ws.Range("K132").Value = 3376.6875
$ws.Range("M132").Value = -846.6875

$ws.Range("H136").Value = 1673.75
$ws.Range("I136").Value = 1673.75
$ws.Range("K136").Value = 5021.25
$ws.Range("M136").Value = -2471.25
